$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.485.04'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '2.281.02'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0898'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('E12').Value = '  -4.67%  '
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.976'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.19'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.57%  '
$ws.Range('D16').Value = '2.628.09'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '2.278.63'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = '42.424.08'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.05%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000104'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '267.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('E24').Value = '  -5.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.20'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.27'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('E33').Value = '  -3.88%  '
$ws.Range('E34').Value = '  -2.11%  '
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('E36').Value = '  -2.82%  '
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  -2.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.76'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  -4.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '108.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +12.69%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.226'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('E46').Value = '  -3.48%  '
$ws.Range('D47').Value = '1.722.02'
$ws.Range('E47').Value = '  +8.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '110.33'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '76.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.73%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.61'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.64%  '
